$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Leaderboard rows (9-46) after the update: new entries inserted and the
# whole board re-sorted by Score descending (ties keep prior relative order).
$data = @(
    ,@('Teacher', 'AAA', 'N/A', 'N/A', 45)
    ,@('Teacher', 'boo', 'Not Applicable', 'Not Applicable', 42)
    ,@('Student', 'Madhu', '5', 'Not Applicable', 38)
    ,@('Teacher', 'rl', 'Not Applicable', 'Not Applicable', 37)
    ,@('Teacher', 'hy', 'Not Applicable', 'Not Applicable', 29)
    ,@('Teacher', 'Abhi', 'Not Applicable', 'Not Applicable', 29)
    ,@('Student', 'Aditya', '7', 'E', 29)
    ,@('Student', '3', '7', 'e', 27)
    ,@('Teacher', 'M', 'Not Applicable', 'Not Applicable', 27)
    ,@('Teacher', 'r', 'Not Applicable', 'Not Applicable', 24)
    ,@('Teacher', 'AAA', 'Not Applicable', 'Not Applicable', 22)
    ,@('Teacher', 'A', 'Not Applicable', 'Not Applicable', 20)
    ,@('Teacher', 'Ab', 'Not Applicable', 'Not Applicable', 19)
    ,@('Teacher', 'drt', 'Not Applicable', 'Not Applicable', 17)
    ,@('Student', 'Aditya', '6', 'B', 15)
    ,@('Teacher', 'rbb', 'Not Applicable', 'Not Applicable', 13)
    ,@('Teacher', 'f', 'Not Applicable', 'Not Applicable', 12)
    ,@('Teacher', 'buck', 'Not Applicable', 'Not Applicable', 12)
    ,@('Teacher', 'MB', 'Not Applicable', 'Not Applicable', 11)
    ,@('Teacher', 'd', 'Not Applicable', 'Not Applicable', 10)
    ,@('Teacher', 'im', 'Not Applicable', 'Not Applicable', 9)
    ,@('Teacher', 'Abhijit', 'Not Applicable', 'Not Applicable', 8)
    ,@('Parent', 'Hello', 'Not Applicable', 'Not Applicable', 7)
    ,@('Teacher', 'hapoochu', 'Not Applicable', 'Not Applicable', 7)
    ,@('Teacher', 'Abh', 'Not Applicable', 'Not Applicable', 7)
    ,@('Parent', 'Aditya', 'Not Applicable', 'Not Applicable', 6)
    ,@('Teacher', 'AAAA', 'Not Applicable', 'Not Applicable', 5)
    ,@('Teacher', 'AAAAA', 'Not Applicable', 'Not Applicable', 5)
    ,@('Teacher', 'AA', 'Not Applicable', 'Not Applicable', 4)
    ,@('Teacher', 'nk', 'Not Applicable', 'Not Applicable', 4)
    ,@('Teacher', 'Montu', 'Not Applicable', 'Not Applicable', 2)
    ,@('Teacher', 'mB', 'Not Applicable', 'Not Applicable', 2)
    ,@('Teacher', 'Mama', 'Not Applicable', 'Not Applicable', 2)
    ,@('Teacher', 't', 'Not Applicable', 'Not Applicable', 1)
    ,@('Student', 'AAAAA', 'LKG', 'E', 1)
    ,@('Teacher', 'B', 'Not Applicable', 'Not Applicable', 0)
    ,@('Teacher', 'AB', 'Not Applicable', 'Not Applicable', 0)
    ,@('Teacher', 'Ad', 'Not Applicable', 'Not Applicable', 0)
)

$startRow = 9
$endRow = $startRow + $data.Count - 1

# Columns A:D (Type/Name/Class/Section) are always text in this sheet, even
# when a value looks numeric (e.g. class "5", "7"); force text formatting
# before writing so Excel does not auto-coerce those cells to numbers.
$ws.Range("A{0}:D{1}" -f $startRow, $endRow).NumberFormat = "@"

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}

